$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - North America
$ws.Range("B2").Value = 126822971
$ws.Range("D2").Value = 1635327
$ws.Range("F2").Value = 122775364
$ws.Range("G2").Value = 4343
$ws.Range("H2").Value = 2412280
$ws.Range("I2").Value = 6014

# Row 3 - Asia
$ws.Range("B3").Value = 217735242
$ws.Range("C3").Value = 1220
$ws.Range("D3").Value = 1547009
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 201469871
$ws.Range("G3").Value = 25229
$ws.Range("H3").Value = 14718362
$ws.Range("I3").Value = 15191

# Row 4 - Europe
$ws.Range("B4").Value = 249540289
$ws.Range("C4").Value = 173
$ws.Range("D4").Value = 2063821
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 245545159
$ws.Range("G4").Value = 3004
$ws.Range("H4").Value = 1931309
$ws.Range("I4").Value = 5640

# Row 5 - South America
$ws.Range("B5").Value = 68766612
$ws.Range("D5").Value = 1356746
$ws.Range("F5").Value = 66479577
$ws.Range("G5").Value = 511
$ws.Range("H5").Value = 930289
$ws.Range("I5").Value = 10104

# Row 6 - Australia/Oceania
$ws.Range("B6").Value = 14458069
$ws.Range("D6").Value = 28511
$ws.Range("F6").Value = 14291244
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 138314

# Row 7 - Africa
$ws.Range("B7").Value = 12825765
$ws.Range("D7").Value = 258782
$ws.Range("F7").Value = 12086419
$ws.Range("H7").Value = 480564
